# Maltaspor.xlsx roster update
# The player/position/team shared-string pool was reordered and rows 2-19
# were re-pointed so each row now shows the correct player alongside their
# real position and current team. Re-apply the final per-row values
# directly; Excel's COM layer takes care of the underlying shared-string
# table bookkeeping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Quentin Grimes", "SG,SF", "Dallas Mavericks"),
    @("Malik Monk", "SG,SF", "Sacramento Kings"),
    @("Brandon Miller", "SG,SF", "Charlotte Hornets"),
    @("Harrison Barnes", "SF,PF", "San Antonio Spurs"),
    @("Nick Richards", "C", "Charlotte Hornets"),
    @("Isaiah Hartenstein", "C", "Oklahoma City Thunder"),
    @("Julius Randle", "PF", "Minnesota Timberwolves"),
    @("Kentavious Caldwell-Pope", "SG,SF", "Orlando Magic"),
    @("Naz Reid", "PF,C", "Minnesota Timberwolves"),
    @("Damian Lillard", "PG", "Milwaukee Bucks"),
    @("Bam Adebayo", "C", "Miami Heat"),
    @("Cade Cunningham", "PG,SG", "Detroit Pistons"),
    @("Derrick White", "PG,SG", "Boston Celtics"),
    @("Anthony Davis", "PF,C", "Los Angeles Lakers"),
    @("Cameron Johnson", "SF,PF", "Brooklyn Nets"),
    @("LaMelo Ball", "PG,SG", "Charlotte Hornets"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("Malcolm Brogdon", "PG,SG", "Washington Wizards")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
